{"js": "// Each table cell holds exactly one arithmetic expression; the new values\n// below replace the old ones in row-major order (20 rows x 5 cols).\nconst newValues = [\n  \"44+48=\",\n  \"67-16=\",\n  \"77-17=\",\n  \"87+1=\",\n  \"38+6=\",\n  \"25+26=\",\n  \"83-58=\",\n  \"18+13=\",\n  \"58-39=\",\n  \"76+10=\",\n  \"84+9=\",\n  \"97-38=\",\n  \"31+16=\",\n  \"95-94=\",\n  \"22+27=\",\n  \"85-8=\",\n  \"45-17=\",\n  \"93-81=\",\n  \"7+78=\",\n  \"83-0=\",\n  \"10+43=\",\n  \"51-22=\",\n  \"15+70=\",\n  \"97-78=\",\n  \"33+54=\",\n  \"7+45=\",\n  \"47-24=\",\n  \"5+23=\",\n  \"29+0=\",\n  \"97-77=\",\n  \"62+35=\",\n  \"63+18=\",\n  \"6+46=\",\n  \"8+91=\",\n  \"43-34=\",\n  \"55+17=\",\n  \"43-12=\",\n  \"9-1=\",\n  \"76-47=\",\n  \"72-45=\",\n  \"66+14=\",\n  \"30-7=\",\n  \"31+5=\",\n  \"42+40=\",\n  \"14+14=\",\n  \"82-4=\",\n  \"29+61=\",\n  \"92-83=\",\n  \"53+32=\",\n  \"65+28=\",\n  \"99-84=\",\n  \"27-15=\",\n  \"91-21=\",\n  \"60-26=\",\n  \"61+7=\",\n  \"91-67=\",\n  \"74+0=\",\n  \"17+75=\",\n  \"28+22=\",\n  \"95-78=\",\n  \"81-36=\",\n  \"37+12=\",\n  \"45+13=\",\n  \"18+55=\",\n  \"18+56=\",\n  \"87-43=\",\n  \"83-66=\",\n  \"23-11=\",\n  \"72+1=\",\n  \"83-43=\",\n  \"40-35=\",\n  \"53-33=\",\n  \"80-1=\",\n  \"88-19=\",\n  \"27+65=\",\n  \"96-34=\",\n  \"1-1=\",\n  \"88-77=\",\n  \"38-2=\",\n  \"45+39=\",\n  \"50+6=\",\n  \"18+75=\",\n  \"35+37=\",\n  \"97-56=\",\n  \"21+52=\",\n  \"0+31=\",\n  \"84-27=\",\n  \"43-31=\",\n  \"80-34=\",\n  \"89-87=\",\n  \"6+10=\",\n  \"53-21=\",\n  \"75-65=\",\n  \"55-31=\",\n  \"21+57=\",\n  \"40+50=\",\n  \"60-26=\",\n  \"84-3=\",\n  \"95-18=\",\n  \"3+41=\"\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\nconst table = tables.items[0];\n\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\n// Load the cells for every row up front.\nconst rowCells = [];\nfor (const row of rows.items) {\n  row.cells.load(\"items\");\n  rowCells.push(row.cells);\n}\nawait context.sync();\n\n// Load the paragraphs for every cell up front.\nconst cellParagraphs = [];\nfor (const cells of rowCells) {\n  for (const cell of cells.items) {\n    cell.body.paragraphs.load(\"items\");\n    cellParagraphs.push(cell.body.paragraphs);\n  }\n}\nawait context.sync();\n\nlet i = 0;\nfor (const paragraphs of cellParagraphs) {\n  const para = paragraphs.items[0];\n  const range = para.getRange();\n  range.insertText(newValues[i], \"Replace\");\n  i++;\n}\nawait context.sync();", "ps1": "$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n# Each table cell holds exactly one arithmetic expression; the new values\n# below replace the old ones in row-major order (20 rows x 5 cols).\n$newValues = @(\n    '44+48=',\n    '67-16=',\n    '77-17=',\n    '87+1=',\n    '38+6=',\n    '25+26=',\n    '83-58=',\n    '18+13=',\n    '58-39=',\n    '76+10=',\n    '84+9=',\n    '97-38=',\n    '31+16=',\n    '95-94=',\n    '22+27=',\n    '85-8=',\n    '45-17=',\n    '93-81=',\n    '7+78=',\n    '83-0=',\n    '10+43=',\n    '51-22=',\n    '15+70=',\n    '97-78=',\n    '33+54=',\n    '7+45=',\n    '47-24=',\n    '5+23=',\n    '29+0=',\n    '97-77=',\n    '62+35=',\n    '63+18=',\n    '6+46=',\n    '8+91=',\n    '43-34=',\n    '55+17=',\n    '43-12=',\n    '9-1=',\n    '76-47=',\n    '72-45=',\n    '66+14=',\n    '30-7=',\n    '31+5=',\n    '42+40=',\n    '14+14=',\n    '82-4=',\n    '29+61=',\n    '92-83=',\n    '53+32=',\n    '65+28=',\n    '99-84=',\n    '27-15=',\n    '91-21=',\n    '60-26=',\n    '61+7=',\n    '91-67=',\n    '74+0=',\n    '17+75=',\n    '28+22=',\n    '95-78=',\n    '81-36=',\n    '37+12=',\n    '45+13=',\n    '18+55=',\n    '18+56=',\n    '87-43=',\n    '83-66=',\n    '23-11=',\n    '72+1=',\n    '83-43=',\n    '40-35=',\n    '53-33=',\n    '80-1=',\n    '88-19=',\n    '27+65=',\n    '96-34=',\n    '1-1=',\n    '88-77=',\n    '38-2=',\n    '45+39=',\n    '50+6=',\n    '18+75=',\n    '35+37=',\n    '97-56=',\n    '21+52=',\n    '0+31=',\n    '84-27=',\n    '43-31=',\n    '80-34=',\n    '89-87=',\n    '6+10=',\n    '53-21=',\n    '75-65=',\n    '55-31=',\n    '21+57=',\n    '40+50=',\n    '60-26=',\n    '84-3=',\n    '95-18=',\n    '3+41='\n)\n\n$cols = $t.Columns.Count\n$i = 0\nforeach ($row in 1..$t.Rows.Count) {\n    foreach ($col in 1..$cols) {\n        $cell = $t.Cell($row, $col)\n        $r = $cell.Range\n        $r.End = $r.End - 2  # drop the trailing cell-mark chars\n        $r.Text = $newValues[$i]\n        $i = $i + 1\n    }\n}\n\nWrite-Output \"updated $i cells\""}
